$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows for "CV + tfidf + ngram(3) + RUS" configuration (Random Under Sampling)
# Columns: A=Model, B=Configuration, C=F1, D=Accuracy, E=Precision, F=Recall
$rusConfig = "CV + tfidf + ngram(3) + RUS"

$data = @(
    @("Logistic Regression", $rusConfig, 88.81, 83.34, 85.39, 92.7),
    @("Multinomial Naive Bayes", $rusConfig, 86.29, 80.37, 86.27, 86.42),
    @("Support Vector Machines", $rusConfig, 87.81, 80.25, 78.53, 99.89),
    @("Decision Tree", $rusConfig, 69.42, 62.01, 83.55, 59.68),
    @("Random Forest", $rusConfig, 77.2, 69.61, 84.2, 71.79)
)

$startRow = 32
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
    $ws.Cells.Item($row, 4).Value = $data[$i][3]
    $ws.Cells.Item($row, 5).Value = $data[$i][4]
    $ws.Cells.Item($row, 6).Value = $data[$i][5]
}

$ws.Range("F37").Select()
